# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Swap the "Periodo Mora" (period) and its corresponding date between the
# two data rows (16 and 17) of the account-statement table:
#   Row16: 1901 / 43355  ->  1812 / 44904
#   Row17: 1812 / 44904  ->  1901 / 43355

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "1812"
$ws.Range("F16").Value = 44904

$ws.Range("E17").Value = "1901"
$ws.Range("F17").Value = 43355
